# Update workbook text to reference the new release:
#   "Coal Mine Boundaries and Methane Sources - version 1.0.0
#    (built on February 03 2026 17.29.55 EST)"
# in place of the old "mines - version 1.0.0 (Feb 3 2026) (built on
# February 03 2026 10.14.00 EST)" build label.

$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Polosukhinskaya Coal Mine, Russia, M1518, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)`""

# --- "Boundaries and methane sources" sheet ---------------------------
# Column S ("build_version") holds the same version label for every
# data row (rows 2 through 17). Row 1 is the header and stays untouched.
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($r = 2; $r -le 17; $r++) {
    $wsData.Range("S$r").Value = $newVersion
}
